$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New TPM-derived row data for rows 2-4 (sending cluster is now "MuSCs" for all three;
# old rows 5-7, the "ECs" sending-cluster rows, are removed).
$rows = @(
    @{ A="MuSCs"; B="Areg"; C="Egfr"; D="ECs";
       E=1; F=0.3333333333333333;
       G=0.014328; H=0.042984; I=1; J=1;
       K=3; L=1;
       M=0.428743; N=1.286229;
       O=0.00412050394863168; P=0.00412050394863168;
       Q=0.006143029704000001; R=0.05528726733600001;
       S=0.00412050394863168; T=0.00412050394863168 },
    @{ A="MuSCs"; B="Areg"; C="Egfr"; D="FAPs";
       E=1; F=0.3333333333333333;
       G=0.014328; H=0.042984; I=1; J=1;
       K=3; L=1;
       M=80.22623699999998; N=240.678711;
       O=0.7710272268990069; P=0.7710272268990069;
       Q=1.149481523736; R=10.345333713624;
       S=0.7710272268990069; T=0.7710272268990069 },
    @{ A="MuSCs"; B="Areg"; C="Egfr"; D="MuSCs";
       E=1; F=0.3333333333333333;
       G=0.014328; H=0.042984; I=1; J=1;
       K=3; L=1;
       M=23.39612766666667; N=70.188383;
       O=0.2248522691523614; P=0.2248522691523614;
       Q=0.335219717208; R=3.016977454872;
       S=0.2248522691523614; T=0.2248522691523614 }
)

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $rowData = $rows[$i]
    foreach ($c in $cols) {
        $ws.Range("$c$r").Value = $rowData[$c]
    }
}

# Remove old rows 5-7 (the former "ECs" sending-cluster rows) entirely.
$ws.Range("A5:T7").EntireRow.Delete()
